$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the whole target range as Text so that numeric-looking
# entries (e.g. "324", "1234") are stored as text strings rather than
# being auto-coerced to numbers - matching the "numberStoredAsText"
# scouting-form data in the source sheet.
$ws.Range("A1:E5").NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "matchNum"
$ws.Range("B1").Value = "TeamNum"
$ws.Range("C1").Value = "climb"
$ws.Range("D1").Value = "throw"
$ws.Range("E1").Value = "additionalNotes"

# Row 2
$ws.Range("A2").Value = "324"
$ws.Range("B2").Value = "1234"
$ws.Range("C2").Value = "No"
$ws.Range("D2").Value = "Yes"
$ws.Range("E2").Value = "49iulfkhdjkhlksajd23"

# Row 3
$ws.Range("A3").Value = "341"
$ws.Range("B3").Value = "2341234"
$ws.Range("C3").Value = "No"
$ws.Range("D3").Value = "Yes"
$ws.Range("E3").Value = "sd32wedaslr2h14lhrkjasgkh"

# Row 4 (new)
$ws.Range("A4").Value = "231"
$ws.Range("B4").Value = "3214"
$ws.Range("C4").Value = "No"
$ws.Range("D4").Value = "Yes"
$ws.Range("E4").Value = "hdfkjhlskdf"

# Row 5 (new)
$ws.Range("A5").Value = "53"
$ws.Range("B5").Value = "4123"
$ws.Range("C5").Value = "Yes"
$ws.Range("D5").Value = "No"
$ws.Range("E5").Value = "421341234"
